$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A30").Value = 29
    $ws.Range("B30").Value = "'2026-02-16"
    $ws.Range("C30").Value = "22:54:59"
    $ws.Range("D30").Value = "base_strategy"
    $ws.Range("E30").Value = "UP"
    $ws.Range("F30").Value = 49.999998
    $ws.Range("G30").Value = ""
    $ws.Range("H30").Value = "OPEN"
    $ws.Range("I30").Value = 0
    $ws.Range("J30").Value = 0
    $ws.Range("K30").Value = 100
    $ws.Range("L30").Value = 0
    $ws.Range("M30").Value = 0
    $ws.Range("N30").Value = 0.6
    $ws.Range("O30").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P30").Value = ""
    $ws.Range("Q30").Value = 0
}
